# Update "想去人数" (interested-attendee counts) figures in column F
# across the 展览, 演出 and 全部类型 sheets, matching the upstream
# site's regenerated data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet
$ws1.Range("F2").Value  = 2434
$ws1.Range("F3").Value  = 675
$ws1.Range("F10").Value = 905
$ws1.Range("F16").Value = 23207
$ws1.Range("F17").Value = 1668
$ws1.Range("F18").Value = 125
$ws1.Range("F20").Value = 19
$ws1.Range("F21").Value = 338
$ws1.Range("F23").Value = 33
$ws1.Range("F26").Value = 24
$ws1.Range("F28").Value = 314

# 演出 sheet
$ws2.Range("F17").Value = 4097

# 全部类型 sheet
$ws4.Range("F4").Value  = 2434
$ws4.Range("F6").Value  = 675
$ws4.Range("F18").Value = 905
$ws4.Range("F23").Value = 23207
$ws4.Range("F29").Value = 1669
$ws4.Range("F30").Value = 125
$ws4.Range("F34").Value = 19
$ws4.Range("F35").Value = 338
$ws4.Range("F37").Value = 33
$ws4.Range("F44").Value = 314
$ws4.Range("F48").Value = 4097
